$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = $true
$ws.Range("B3").Value = $true
$ws.Range("B4").Value = $true
